# Add 10 new rows of Mac-Address/device data to the
# master-reg_center_device_h sheet (rows 147-156), continuing the
# existing pattern of regcntr_id/device_id/lang_code/is_active/cr_by/
# cr_dtimes/eff_dtimes values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 146; new rows start at 147.
$startRow = 147
$startDevice = 3000166
$rowCount = 10

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 10001
    $ws.Cells.Item($r, 2).Value = $startDevice + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

$lastRow = $startRow + $rowCount - 1

# Update the active selection to match where the user ended up after
# typing the new data.
$ws.Range("E155").Select()
